$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row with the approximate calibrated R0 for each stage (string
# literal written first so it lands at the same shared-string index the
# saved workbook uses).
$ws.Range("A14").Value = "Approx. Calibrate R0"

# Insert a new column before old column D ("Stage 3"), shifting old D->E, E->F, F->G.
# This makes room for a new "Stage 2b" column, fixing "stage 2 skipping".
$ws.Columns("D:D").Insert()

$ws.Range("D1").Value = "Stage 2b"

# Fill in the new column D values for each parameter row.
$ws.Range("D2").Value = 6.5
$ws.Range("D3").Value = 70
$ws.Range("D4").Value = 70
$ws.Range("D5").Value = 60
$ws.Range("D6").Value = 45
$ws.Range("D7").Value = $true
$ws.Range("D8").Value = 55
$ws.Range("D9").Value = 6
$ws.Range("D10").Value = 0.13
$ws.Range("D11").Value = 0.0064
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = "NA"

# Complete row 14 with the per-stage R0 values.
$ws.Range("B14").Value = 1.82
$ws.Range("C14").Value = 1.35
$ws.Range("D14").Value = 1.1
$ws.Range("E14").Value = 0.89
$ws.Range("F14").Value = 0.77

# Apply a thin box border to the whole parameter table (A1:F14).
$borderRange = $ws.Range("A1:F14")
$borderRange.Borders.LineStyle = 1
$borderRange.Borders.Weight = 2

# Cosmetic view adjustments matching the saved workbook state.
$ws.Range("L20").Select() | Out-Null
